$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# New ordering / labels / values for rows 2-34 (row 1 "Danh mục lương" stays as-is)
$rows = @(
    @("Tổng công tại CẦN THƠ", 28),
    @("Lương cơ bản tại CẦN THƠ", 3000000),
    @("Chiết khấu sale chính tại CẦN THƠ", 260000),
    @("Chiết khấu sale phụ tại CẦN THƠ", 0),
    @("Đơn 1 bác sĩ tại CẦN THƠ", 0),
    @("Đơn 2 bác sĩ tại CẦN THƠ", 0),
    @("Công phụ phẫu 1 tại CẦN THƠ", 50000),
    @("Công phụ phẫu 2 tại CẦN THƠ", 0),
    @("Ứng lương tại CẦN THƠ", -3850000),
    @("Tổng công tại LONG XUYÊN", 0),
    @("Lương công tác tại LONG XUYÊN", 0),
    @("Lương cơ bản tại LONG XUYÊN", 2000000),
    @("Chiết khấu sale chính tại LONG XUYÊN", 0),
    @("Chiết khấu sale phụ tại LONG XUYÊN", 0),
    @("Đơn 1 bác sĩ tại LONG XUYÊN", 0),
    @("Đơn 2 bác sĩ tại LONG XUYÊN", 0),
    @("Công phụ phẫu 1 tại LONG XUYÊN", 0),
    @("Công phụ phẫu 2 tại LONG XUYÊN", 0),
    @("Ứng lương tại LONG XUYÊN", -0),
    @("Tổng công tại SÓC TRĂNG", 0),
    @("Lương công tác tại SÓC TRĂNG", 0),
    @("Lương cơ bản tại SÓC TRĂNG", 3000000),
    @("Chiết khấu sale chính tại SÓC TRĂNG", 0),
    @("Chiết khấu sale phụ tại SÓC TRĂNG", 0),
    @("Đơn 1 bác sĩ tại SÓC TRĂNG", 0),
    @("Đơn 2 bác sĩ tại SÓC TRĂNG", 0),
    @("Công phụ phẫu 1 tại SÓC TRĂNG", 0),
    @("Công phụ phẫu 2 tại SÓC TRĂNG", 0),
    @("Ứng lương tại SÓC TRĂNG", -0),
    @("Tổng lương tại CẦN THƠ", -539999.9999999995),
    @("Tổng lương tại LONG XUYÊN", 2000000),
    @("Tổng lương tại SÓC TRĂNG", 3000000),
    @("Tổng lương", 4460000.000000001)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}
